$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the four section-header labels in column A. These are plain text
# (shared-string) edits only -- the row/column layout, styles and every
# other cell in the sheet stay exactly as they were.
$ws.Range("A2").Value  = "market_paries"
$ws.Range("A6").Value  = "net_areas"
$ws.Range("A12").Value = "grid_points"
$ws.Range("A18").Value = "supply_contracts"

# Move the active selection to D21 (was D15), matching where the author was
# last working in the sheet.
$ws.Range("D21").Select()
